$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Alessio Zandonai" team row at the bottom of the table (row 31)
$ws.Range("A31").Value = "Alessio Zandonai"
$ws.Range("B31").Value = "Daniele  Dalbosco | iMontagna"
$ws.Range("C31").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("D31").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E31").Value = "Luca Frasca | Clitoriders"
$ws.Range("F31").Value = "Davide  Bazzano  | iMontagna"
